{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paras = body.paragraphs.items;\n\n// Original paragraph order (0-indexed):\n// 0: \"Nombre: Alexis Bonilla\"\n// 1: \"Curso: Programaci\u00f3n web avanzada.\"\n// 2: (empty) -- stays untouched\n// 3: \"Taller 1: HTML + CSS +JS \"\n// 4: (empty)\n// 5: \"Definitivamente este taller me sirvi\u00f3...\"\n// 6: \"Al intentar replicar el template...\"\n// 7: \"Excelente ejercicio para aprender...\" (multi-run)\n// 8: (empty)\n// 9: (empty)\n// 10: \"Link del repositorio: \" + url\n\nconst pName = paras[0];\nconst pCourse = paras[1];\nconst pWorkshop = paras[3];\nconst pEmptyAfterWorkshop = paras[4];\nconst pDef = paras[5];\nconst pAl = paras[6];\nconst pExcelente = paras[7];\nconst pEmpty8 = paras[8];\nconst pEmpty9 = paras[9];\nconst pLink = paras[10];\n\n// Remove the paragraphs that collapse away entirely.\npCourse.delete();\npEmptyAfterWorkshop.delete();\npDef.delete();\npAl.delete();\npExcelente.delete();\npEmpty8.delete();\npEmpty9.delete();\npLink.delete();\n\n// Rewrite the two paragraphs that survive (in place) with their new text.\n// clear() first so the new text doesn't inherit the old run's\n// xml:space=\"preserve\" (the old \"Taller 1: ... \" text ended with a space).\npName.clear();\npName.insertText(\"Alexis Bonilla- Resumen de aprendizaje.\", Word.InsertLocation.replace);\n\npWorkshop.clear();\npWorkshop.insertText(\n  \"Para lograr hacer la p\u00e1gina con Bootstrap, tuve que aprender a utilizar la web de Bootstrap 4 para averiguar las diferentes funcionalidades de los componentes que utilic\u00e9, definitivamente el reto m\u00e1s importante a la hora de realizar la p\u00e1gina fue hacer la barra de navegaci\u00f3n vertical, ya que en clase no la hab\u00edamos visto. Aprend\u00ed mucho acerca de inspeccionar una p\u00e1gina ya hecha y c\u00f3mo entender el c\u00f3digo de otros. Este ejercicio es de gran ayuda para mi futuro profesional, ya que entre las cosas m\u00e1s importantes que hay en nuestra carrera es aprender a aprender.\",\n  Word.InsertLocation.replace\n);\n\n// Insert a new leading empty paragraph before the (now-renamed) first paragraph.\npName.insertParagraph(\"\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Original paragraph layout (1-indexed):\n#  1: \"Nombre: Alexis Bonilla\"\n#  2: \"Curso: Programaci\u00f3n web avanzada.\"\n#  3: (empty) -- stays untouched\n#  4: \"Taller 1: HTML + CSS +JS \"\n#  5: (empty)\n#  6: \"Definitivamente este taller me sirvi\u00f3 bastante...\"\n#  7: \"Al intentar replicar el template del ejercicio...\"\n#  8: \"Excelente ejercicio para aprender todo...\"\n#  9: (empty)\n# 10: (empty)\n# 11: \"Link del repositorio: \" + url\n\n# Remove paragraphs 5 through 11 in one pass (a single contiguous Range\n# delete avoids the collection re-indexing entirely).\n$pFrom = $d.Paragraphs.Item(5)\n$pTo = $d.Paragraphs.Item(11)\n$d.Range($pFrom.Range.Start, $pTo.Range.End).Delete()\n\n# Remove paragraph 2 (\"Curso: ...\").\n$d.Paragraphs.Item(2).Range.Delete()\n\n# Rewrite the two paragraphs that survive (in place) with their new text,\n# using Find/Replace (wdReplaceAll) so the surviving run stays clean instead\n# of inheriting the old run's xml:space=\"preserve\" (the old \"Taller 1: ... \"\n# text ended with a trailing space).\n$d.Content.Find.Execute(\n    \"Nombre: Alexis Bonilla\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Alexis Bonilla- Resumen de aprendizaje.\", 2) | Out-Null\n\n$d.Content.Find.Execute(\n    \"Taller 1: HTML + CSS +JS \", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Para lograr hacer la p\u00e1gina con Bootstrap, tuve que aprender a utilizar la web de Bootstrap 4 para averiguar las diferentes funcionalidades de los componentes que utilic\u00e9, definitivamente el reto m\u00e1s importante a la hora de realizar la p\u00e1gina fue hacer la barra de navegaci\u00f3n vertical, ya que en clase no la hab\u00edamos visto. Aprend\u00ed mucho acerca de inspeccionar una p\u00e1gina ya hecha y c\u00f3mo entender el c\u00f3digo de otros. Este ejercicio es de gran ayuda para mi futuro profesional, ya que entre las cosas m\u00e1s importantes que hay en nuestra carrera es aprender a aprender.\", 2) | Out-Null\n\n# Insert a new leading empty paragraph before the (now-renamed) first paragraph.\n$d.Paragraphs.Item(1).Range.InsertParagraphBefore()\n"}
